$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the JSON text in-place first (D2 still holds the sole reference to
# the old shared string) so it keeps the same shared-string slot, just with
# "type" renamed to "target".
$ws.Range("D2").Value = "{`n    ""value"":""echo `${name}_`${caseId}"",`n    ""target"":""cmd""`n}"

# Insert a new column before column E, copying column D's formatting/width
# so the new E (run-target) column looks like the existing D (run) column.
$ws.Columns.Item(4).Copy()
$ws.Columns.Item(5).Insert()
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(4).ColumnWidth

# --- Row 2: split the old "run" JSON args into D2 ("cmd") / E2 (json) ---
# D2 now holds the short literal that used to be the "type" value; E2 (the
# copy) keeps the updated JSON blob.
$ws.Range("D2").Value = "cmd"

# --- Row 3: give D3 a new value + its own (bordered/wrapped) style ---
$ws.Range("D3").Value = "echo `${name}"
$ws.Range("D3").Borders.LineStyle = 1
$ws.Range("D3").HorizontalAlignment = -4131
$ws.Range("D3").VerticalAlignment = -4108
$ws.Range("D3").WrapText = $true

# Update the active selection to match the new layout.
$ws.Range("E5").Select()
